$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(11, 1).Value = '2024-10-08 21:15:44'
$ws.Cells.Item(11, 2).Value = 'get_price'
$ws.Cells.Item(11, 3).Value = 'https://example.com/product'
$ws.Cells.Item(11, 4).Value = 'Error fetching price: ''NoneType'' object is not subscriptable'
$ws.Cells.Item(11, 5).NumberFormat = '@'
$ws.Cells.Item(11, 5).Value = '2024-10-08'
$ws.Cells.Item(11, 5).Style = 'Normal'
$ws.Cells.Item(11, 6).Value = '21:15:44'

$ws.Cells.Item(12, 1).Value = '2024-10-08 21:15:47'
$ws.Cells.Item(12, 2).Value = 'get_price'
$ws.Cells.Item(12, 3).Value = 'https://example.com/product'
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '100.00'
$ws.Cells.Item(12, 4).Style = 'Normal'
$ws.Cells.Item(12, 5).NumberFormat = '@'
$ws.Cells.Item(12, 5).Value = '2024-10-08'
$ws.Cells.Item(12, 5).Style = 'Normal'
$ws.Cells.Item(12, 6).Value = '21:15:47'

$ws.Cells.Item(13, 1).Value = '2024-10-08 21:15:52'
$ws.Cells.Item(13, 2).Value = 'get_price'
$ws.Cells.Item(13, 3).Value = 'https://example.com/product'
$ws.Cells.Item(13, 4).Value = 'Error fetching price: ''NoneType'' object is not subscriptable'
$ws.Cells.Item(13, 5).NumberFormat = '@'
$ws.Cells.Item(13, 5).Value = '2024-10-08'
$ws.Cells.Item(13, 5).Style = 'Normal'
$ws.Cells.Item(13, 6).Value = '21:15:52'

$ws.Cells.Item(14, 1).Value = '2024-10-08 21:18:33'
$ws.Cells.Item(14, 2).Value = 'get_price'
$ws.Cells.Item(14, 3).Value = 'https://example.com/product'
$ws.Cells.Item(14, 4).Value = 'Error fetching price: ''NoneType'' object is not subscriptable'
$ws.Cells.Item(14, 5).NumberFormat = '@'
$ws.Cells.Item(14, 5).Value = '2024-10-08'
$ws.Cells.Item(14, 5).Style = 'Normal'
$ws.Cells.Item(14, 6).Value = '21:18:33'

$ws.Cells.Item(15, 1).Value = '2024-10-08 21:18:35'
$ws.Cells.Item(15, 2).Value = 'get_price'
$ws.Cells.Item(15, 3).Value = 'https://example.com/product'
$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '100.00'
$ws.Cells.Item(15, 4).Style = 'Normal'
$ws.Cells.Item(15, 5).NumberFormat = '@'
$ws.Cells.Item(15, 5).Value = '2024-10-08'
$ws.Cells.Item(15, 5).Style = 'Normal'
$ws.Cells.Item(15, 6).Value = '21:18:35'

$ws.Cells.Item(16, 1).Value = '2024-10-08 21:18:41'
$ws.Cells.Item(16, 2).Value = 'get_price'
$ws.Cells.Item(16, 3).Value = 'https://example.com/product'
$ws.Cells.Item(16, 4).Value = 'Error fetching price: ''NoneType'' object is not subscriptable'
$ws.Cells.Item(16, 5).NumberFormat = '@'
$ws.Cells.Item(16, 5).Value = '2024-10-08'
$ws.Cells.Item(16, 5).Style = 'Normal'
$ws.Cells.Item(16, 6).Value = '21:18:41'

$ws.Cells.Item(17, 1).Value = '2024-10-08 21:24:56'
$ws.Cells.Item(17, 2).Value = 'get_price'
$ws.Cells.Item(17, 3).Value = 'https://example.com/product'
$ws.Cells.Item(17, 4).Value = 'Error fetching price: ''NoneType'' object is not subscriptable'
$ws.Cells.Item(17, 5).NumberFormat = '@'
$ws.Cells.Item(17, 5).Value = '2024-10-08'
$ws.Cells.Item(17, 5).Style = 'Normal'
$ws.Cells.Item(17, 6).Value = '21:24:56'

$ws.Cells.Item(18, 1).Value = '2024-10-08 21:24:58'
$ws.Cells.Item(18, 2).Value = 'get_price'
$ws.Cells.Item(18, 3).Value = 'https://example.com/product'
$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '100.00'
$ws.Cells.Item(18, 4).Style = 'Normal'
$ws.Cells.Item(18, 5).NumberFormat = '@'
$ws.Cells.Item(18, 5).Value = '2024-10-08'
$ws.Cells.Item(18, 5).Style = 'Normal'
$ws.Cells.Item(18, 6).Value = '21:24:58'

$ws.Cells.Item(19, 1).Value = '2024-10-08 21:25:03'
$ws.Cells.Item(19, 2).Value = 'get_price'
$ws.Cells.Item(19, 3).Value = 'https://example.com/product'
$ws.Cells.Item(19, 4).Value = 'Error fetching price: ''NoneType'' object is not subscriptable'
$ws.Cells.Item(19, 5).NumberFormat = '@'
$ws.Cells.Item(19, 5).Value = '2024-10-08'
$ws.Cells.Item(19, 5).Style = 'Normal'
$ws.Cells.Item(19, 6).Value = '21:25:03'

